# Auto commit: 2024-10-16 19:27:26
# Adds a new post row (row 12) to the "posts" sheet, extends the C2:C12
# data-validation range, adjusts the sheet view, and applies the banded
# border/fill formatting that was added to the table.

$wb = $excel.ActiveWorkbook

$posts = $wb.Worksheets.Item("posts")
$topics = $wb.Worksheets.Item("topics")

# ---------------------------------------------------------------------
# 1. New row of data (row 12)
# ---------------------------------------------------------------------
$posts.Range("A12").Value = 11
$posts.Range("B12").Value = "Quickly Integrate Tailwind CSS into an Existing Next.js App in Just 4 Steps"
$posts.Range("C12").Value = "programming"
$posts.Range("D12").Value = "minimalistic-cartoon-style-laptop-nextjs-app-tailwindcss-color-swatches-code-snippets-black-outline"
$posts.Range("E12").Value = "#A2F9E9"
$posts.Range("F12").Value = "Recent versions of Next.js come packaged with Tailwind CSS, which we can enable during project creation. If you missed installing it out of the box, this guide will walk you through four simple steps to install and configure Tailwind CSS. Even though this post is tailored for Next.js, a similar method can be used to set up Tailwind in a React.js application as well."

$posts.Range("G12").Formula = "=LOWER(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(B12, "" "", ""-""), "":"", ""-""), "";"", ""-""), "","", ""-""), ""."", ""-""), ""/"", ""-""))"
$posts.Range("H12").Formula = "=SUBSTITUTE(SUBSTITUTE(D12, ""-"", "" ""), ""."", "" "")"

$posts.Range("I12").Value = $true
$posts.Range("J12").Value = $false
$posts.Range("K12").Value = $false
$posts.Range("L12").Value = 45581
$posts.Range("M12").Value = 45581

# ---------------------------------------------------------------------
# 2. Data validation on TOPIC now covers the new row too
# ---------------------------------------------------------------------
$posts.Range("C2:C12").Validation.Delete()
$posts.Range("C2:C12").Validation.Add(3, 1, 1, "=topics!`$C`$2:`$C`$1000")

# ---------------------------------------------------------------------
# 3. Formatting: thin border around every cell in the table, plus the
#    banded fill colours that were introduced together with the new row
# ---------------------------------------------------------------------
$allTable = $posts.Range("A1:M12")
$allTable.Borders.LineStyle = 1
$allTable.Borders.Weight = 2

$posts.Range("A2:M11").Interior.Color = 13431551
$posts.Range("I12:K12").Interior.Color = 16777215

# ---------------------------------------------------------------------
# 4. View state: scroll + selection
# ---------------------------------------------------------------------
$posts.Application.ActiveWindow.ScrollRow = 2
$posts.Range("B11").Select()
